$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the end-time entry for the row-29 work block; dependent formulas
# (duration in D29, and the cumulative totals in F29:F33 / G29:G33 that
# chain off of it) recalculate automatically.
$ws.Range("C29").Value = 0.90625

# Move the active cell/selection to C30, matching the saved view state.
[void]$ws.Range("C30").Select()
